$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-34 down to 18-35
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new record's data
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Macroferia Regional de Talca"
$ws.Range("C17").Value = "Maule"
$ws.Range("D17").Value = 44566
$ws.Range("E17").Value = 7
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100103
$ws.Range("H17").Value = "Frutos de hueso (carozo)"
$ws.Range("I17").Value = 100103003
$ws.Range("J17").Value = "Damasco"
$ws.Range("K17").Value = "Modesto"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 10000
$ws.Range("O17").Value = 10000
$ws.Range("P17").Value = 10000
$ws.Range("Q17").Value = "$/caja 10 kilos"
$ws.Range("R17").Value = "Región de O'Higgins"
$ws.Range("S17").Value = 1000
$ws.Range("T17").Value = 10
